$d = $word.ActiveDocument

function Replace-Exact($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Exact "617×4=" "828×9="
Replace-Exact "657×3=" "602×2="
Replace-Exact "434×8=" "261×3="
Replace-Exact "399×7=" "426×3="
Replace-Exact "573×3=" "613×2="
Replace-Exact "487×2=" "345×4="
Replace-Exact "218×5=" "820×6="
Replace-Exact "429×6=" "981×2="
Replace-Exact "585×2=" "316×2="
Replace-Exact "132×9=" "459×5="
Replace-Exact "448×8=" "163×2="
Replace-Exact "369×2=" "516×9="
Replace-Exact "972×2=" "607×4="
Replace-Exact "435×4=" "904×5="
Replace-Exact "764×8=" "549×9="
Replace-Exact "574×9=" "843×6="
Replace-Exact "276×6=" "584×6="
Replace-Exact "982×3=" "419×2="
Replace-Exact "312×2=" "472×3="
Replace-Exact "604×7=" "540×2="
Replace-Exact "293×3=" "265×7="
Replace-Exact "965×2=" "852×8="
Replace-Exact "185×2=" "743×9="
Replace-Exact "692×7=" "390×3="
Replace-Exact "311×6=" "807×5="

Write-Output "Done"
